$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2:E2").NumberFormat = "@"
$ws.Range("D2").Value = "320.14"
$ws.Range("E2").Value = "-3.62%"
$ws.Range("D2:E2").Style = "Normal"

# Row 3
$ws.Range("D3:E3").NumberFormat = "@"
$ws.Range("D3").Value = "42.66"
$ws.Range("E3").Value = "-7.03%"
$ws.Range("D3:E3").Style = "Normal"

# Row 4
$ws.Range("D4:E4").NumberFormat = "@"
$ws.Range("D4").Value = "5.239"
$ws.Range("E4").Value = "-7.66%"
$ws.Range("D4:E4").Style = "Normal"

# Row 5
$ws.Range("D5:E5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08149"
$ws.Range("E5").Value = "-2.60%"
$ws.Range("D5:E5").Style = "Normal"

# Row 6
$ws.Range("D6:E6").NumberFormat = "@"
$ws.Range("D6").Value = "4.332"
$ws.Range("E6").Value = "-3.32%"
$ws.Range("D6:E6").Style = "Normal"

# Row 7
$ws.Range("D7:E7").NumberFormat = "@"
$ws.Range("D7").Value = "1.751"
$ws.Range("E7").Value = "-14.16%"
$ws.Range("D7:E7").Style = "Normal"

# Row 8
$ws.Range("D8:E8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9508"
$ws.Range("E8").Value = "-3.67%"
$ws.Range("D8:E8").Style = "Normal"

# Row 9
$ws.Range("D9:E9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1111"
$ws.Range("E9").Value = "-4.40%"
$ws.Range("D9:E9").Style = "Normal"

# Row 10
$ws.Range("D10:E10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1850"
$ws.Range("E10").Value = "-4.79%"
$ws.Range("D10:E10").Style = "Normal"

# Row 11
$ws.Range("B11").Value = 'BitrueCoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D11:E11").NumberFormat = "@"
$ws.Range("D11").Value = "0.04668"
$ws.Range("E11").Value = "0.04%"
$ws.Range("D11:E11").Style = "Normal"

# Row 12
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12:E12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09311"
$ws.Range("E12").Value = "-7.90%"
$ws.Range("D12:E12").Style = "Normal"

# Row 13
$ws.Range("D13:E13").NumberFormat = "@"
$ws.Range("D13").Value = "7.425"
$ws.Range("E13").Value = "-28.84%"
$ws.Range("D13:E13").Style = "Normal"

# Row 14
$ws.Range("D14:E14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1061"
$ws.Range("E14").Value = "0.24%"
$ws.Range("D14:E14").Style = "Normal"

# Row 15
$ws.Range("D15:E15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001282"
$ws.Range("E15").Value = "-0.25%"
$ws.Range("D15:E15").Style = "Normal"

# Row 16
$ws.Range("B16").Value = 'CoinExToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D16:E16").NumberFormat = "@"
$ws.Range("D16").Value = "0.04191"
$ws.Range("E16").Value = "-0.29%"
$ws.Range("D16:E16").Style = "Normal"

# Row 17
$ws.Range("B17").Value = 'TigerCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D17:E17").NumberFormat = "@"
$ws.Range("D17").Value = "0.005918"
$ws.Range("E17").Value = "-3.53%"
$ws.Range("D17:E17").Style = "Normal"

# Row 18
$ws.Range("B18").Value = 'LEO'
$ws.Range("C18").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D18:E18").NumberFormat = "@"
$ws.Range("D18").Value = "3.365"
$ws.Range("E18").Value = "-0.13%"
$ws.Range("D18:E18").Style = "Normal"

# Row 19
$ws.Range("B19").Value = 'BTSEToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D19:E19").NumberFormat = "@"
$ws.Range("D19").Value = "2.538"
$ws.Range("E19").Value = "-1.36%"
$ws.Range("D19:E19").Style = "Normal"

# Row 20
$ws.Range("B20").Value = 'BitpandaEcosystemToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D20:E20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3363"
$ws.Range("E20").Value = "-0.04%"
$ws.Range("D20:E20").Style = "Normal"

# Row 21
$ws.Range("B21").Value = 'ProBitToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D21:E21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1393"
$ws.Range("E21").Value = "-0.52%"
$ws.Range("D21:E21").Style = "Normal"

# Row 22
$ws.Range("B22").Value = 'ZBToken'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("D22:E22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2686"
$ws.Range("E22").Value = "1.47%"
$ws.Range("D22:E22").Style = "Normal"

# Row 23
$ws.Range("D23:E23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001253"
$ws.Range("E23").Value = "-4.25%"
$ws.Range("D23:E23").Style = "Normal"

# Row 24
$ws.Range("D24:E24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004315"
$ws.Range("E24").Value = "-7.40%"
$ws.Range("D24:E24").Style = "Normal"

# Row 25
$ws.Range("D25:E25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001303"
$ws.Range("E25").Value = "1.70%"
$ws.Range("D25:E25").Style = "Normal"

# Row 26
$ws.Range("D26:E26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0002991"
$ws.Range("E26").Value = "-20.08%"
$ws.Range("D26:E26").Style = "Normal"

# Row 38
$ws.Range("D38:E38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02597"
$ws.Range("E38").Value = "-7.02%"
$ws.Range("D38:E38").Style = "Normal"

# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-5.37%"
$ws.Range("E39").Style = "Normal"

# Row 40
$ws.Range("D40:E40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007817"
$ws.Range("E40").Value = "1.08%"
$ws.Range("D40:E40").Style = "Normal"

# Row 41
$ws.Range("D41:E41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1390"
$ws.Range("E41").Value = "-3.30%"
$ws.Range("D41:E41").Style = "Normal"

# Row 42
$ws.Range("D42:E42").NumberFormat = "@"
$ws.Range("D42").Value = "0.006611"
$ws.Range("E42").Value = "-9.22%"
$ws.Range("D42:E42").Style = "Normal"

# Row 43
$ws.Range("D43:E43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002126"
$ws.Range("E43").Value = "7.64%"
$ws.Range("D43:E43").Style = "Normal"

# Row 44
$ws.Range("D44:E44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008479"
$ws.Range("E44").Value = "-7.54%"
$ws.Range("D44:E44").Style = "Normal"

# Row 45
$ws.Range("D45:E45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3464"
$ws.Range("E45").Value = "-1.11%"
$ws.Range("D45:E45").Style = "Normal"

# Row 46
$ws.Range("D46:E46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006982"
$ws.Range("E46").Value = "-2.85%"
$ws.Range("D46:E46").Style = "Normal"

# Row 47
$ws.Range("D47:E47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000752"
$ws.Range("E47").Value = "0.15%"
$ws.Range("D47:E47").Style = "Normal"

# Row 48
$ws.Range("D48:E48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003490"
$ws.Range("E48").Value = "-0.40%"
$ws.Range("D48:E48").Style = "Normal"

# Row 49
$ws.Range("D49:E49").NumberFormat = "@"
$ws.Range("D49").Value = "0.003545"
$ws.Range("E49").Value = "1.25%"
$ws.Range("D49:E49").Style = "Normal"

# Row 50
$ws.Range("D50:E50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002105"
$ws.Range("E50").Value = "0.15%"
$ws.Range("D50:E50").Style = "Normal"

# Row 51
$ws.Range("D51:E51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002005"
$ws.Range("E51").Value = "0.15%"
$ws.Range("D51:E51").Style = "Normal"
